$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ==========================================================================
# Phase 1: tag every text run whose content needs to move, with a unique
# placeholder token, so later replacements cannot accidentally collide with
# text that has not been processed yet.
# ==========================================================================

Replace-Text "Prover conhecimento e ferramentas para análise da sustentabilidade de cadeias produtivas, desenvolvendo um entendimento sobre como usar as decisões de engenharia para melhorar a performance ambiental, social e econômica." "@@P_OBJ_PT@@"
Replace-Text "Provide knowledge and tools for analyzing the sustainability of production chains, developing an understanding of how to use engineering decisions to improve environmental, social and economic performance." "@@P_OBJ_EN@@"
Replace-Text "3295113 - José Eduardo Holler Branco" "@@P_DOC1@@"
Replace-Text "5840535 - Messias Borges Silva" "@@P_DOC2@@"
Replace-Text "Planejamento de cadeias de suprimentos sustentáveis." "@@P_RES_PT@@"
Replace-Text "Planning sustainable supply chains." "@@P_RES_EN@@"
Replace-Text "i) Sustentabilidade da cadeia de suprimentos; ii) Planejamento focado na redução do consumo energético e de combustíveis fósseis; iii) Combustíveis renováveis e de baixa emissão de gases do efeito estufa; iv) O mercado de carbono; v) Programas de certificação e vi) Economia Circular e Logística Reversa." "@@P_PROG_PT@@"
Replace-Text "Aulas Expositivas; trabalhos e seminários" "@@P_METODO_V@@"
Replace-Text "Média das atividades avaliativas." "@@P_CRIT_V@@"
Replace-Text "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação." "@@P_NORMA_V@@"

Write-Output "=== AFTER PHASE 1 (placeholders) ==="
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output "$i : $($d.Paragraphs($i).Range.Text)"
}

# ==========================================================================
# Phase 2: move the Bibliografia run's content (8 entries, currently its own
# paragraph) into the Avaliação paragraph's "Critério" value run, leaving
# only "5840535 - Messias Borges Silva" behind in the old Bibliografia spot.
# ==========================================================================

# First, tag the *original* bibliography paragraph's whole content with a
# placeholder (using a wildcard match) so it cannot be confused with the
# copy we are about to place in the Avaliação paragraph.
$d.Content.Find.Execute("BOWERSOX*São Paulo, 2002.", $true, $false, $true, $false, $false, $true, 1, $false, "@@P_BIB_OLD@@", 2) | Out-Null

Replace-Text "@@P_CRIT_V@@" "BOWERSOX, D. J.; CLOSS, D. J.; COOPER, M. B.; BOWERSOX, J. C. Gestão Logística da Cadeia de Suprimentos. [s.l.] AMGH, 2013. 472 p.`v`vBARBIERI, J. C. Gestão Ambiental Empresarial: conceitos, modelos e instrumentos. Editora Saraiva, 2004.`v`vALLEN, D.T.; SHONNARD, D. R., Sustainable Engineering: concepts, design and case studies, Prentice Hall, 2015. `v`vAKKUCUK, U. Handbook of Research on Sustainable Supply Chain Management for the Global Economy. [s.l.] IGI Global, 2020. 409 p.`v`vBOUCHERY, Y.; CORBETT, C. J.; FRANSOO, J. C.; TAN, T. (ed.). Sustainable Supply Chains. Cham: Springer International Publishing, 2017. v. 4. 130 p.`v`vSCHMIDT, M.; GIOVANNUCCI, D.; PALEKHOV, D.; HANSMANN, B. (ed.). Sustainable Global Value Chains. Cham: Springer International Publishing, 2019. v. 2. 304 p.`v`vLAVE, L. B.; HENDRICKSON, C. T. Environmental Life Cycle Assessment of Goods and Services, Editora John Hopkins, 2006.`v`vLEITE, P. R. Logística Reversa - Meio Ambiente e Competitividade, Editora Prentice Hall: São Paulo, 2002."

Write-Output "=== AFTER PHASE 2 (bibliography moved) ==="
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output "$i : $($d.Paragraphs($i).Range.Text)"
}

# ==========================================================================
# Phase 3: resolve every remaining placeholder to its final text.
# ==========================================================================

Replace-Text "@@P_OBJ_PT@@"   "Planejamento de cadeias de suprimentos sustentáveis."
Replace-Text "@@P_OBJ_EN@@"   "Planning sustainable supply chains."
Replace-Text "@@P_DOC1@@"     "Prover conhecimento e ferramentas para análise da sustentabilidade de cadeias produtivas, desenvolvendo um entendimento sobre como usar as decisões de engenharia para melhorar a performance ambiental, social e econômica."
Replace-Text "@@P_DOC2@@"     "i) Sustentabilidade da cadeia de suprimentos; ii) Planejamento focado na redução do consumo energético e de combustíveis fósseis; iii) Combustíveis renováveis e de baixa emissão de gases do efeito estufa; iv) O mercado de carbono; v) Programas de certificação e vi) Economia Circular e Logística Reversa."
Replace-Text "@@P_RES_PT@@"   "Aulas Expositivas; trabalhos e seminários"
Replace-Text "@@P_RES_EN@@"   "Provide knowledge and tools for analyzing the sustainability of production chains, developing an understanding of how to use engineering decisions to improve environmental, social and economic performance."
Replace-Text "@@P_PROG_PT@@"  "Média das atividades avaliativas."
Replace-Text "@@P_METODO_V@@" "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
Replace-Text "@@P_NORMA_V@@"  "3295113 - José Eduardo Holler Branco"
Replace-Text "@@P_BIB_OLD@@"  "5840535 - Messias Borges Silva"

Write-Output "=== FINAL ==="
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output "$i : $($d.Paragraphs($i).Range.Text)"
}
